# Konduga.xlsx edit: enhanced fields autocomplete with templates
#
# 1) Rows 2-25: the "flag" column (K) currently stores the text "False" as a
#    string. Convert these to real Boolean FALSE values.
# 2) Rows 26-31: refresh the sample/template annotation data (columns
#    B..F, J) and likewise convert column K to a real Boolean FALSE.
# 3) Rows 32-36 (old rows 30-34) are removed entirely, shrinking the sheet
#    to A1:K31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 2-25: flag column -> Boolean FALSE -----------------------------
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 11).Value = $false
}

# --- 2) Rows 26-31: new template content ------------------------------------

# Row 26
$ws.Range("B26").Value = "1:1"
$ws.Range("C26").Value = "لِلهِ"
$ws.Range("D26").Value = "asa"
$ws.Range("E26").Value = "s"
$ws.Range("J26").Value = "12 - 17"
$ws.Range("K26").Value = $false

# Row 27
$ws.Range("B27").Value = "1:1"
$ws.Range("C27").Value = "لِلهِ"
$ws.Range("D27").Value = "asa"
$ws.Range("E27").Value = "s"
$ws.Range("J27").Value = "12 - 17"
$ws.Range("K27").Value = $false

# Row 28
$ws.Range("B28").Value = "1:1"
$ws.Range("C28").Value = "اِ۬لْحَمْدُ"
$ws.Range("D28").Value = "aqaq"
$ws.Range("E28").Value = "Ar"
$ws.Range("F28").Value = "qwqw"
$ws.Range("J28").Value = "0 - 11"
$ws.Range("K28").Value = $false

# Row 29
$ws.Range("B29").Value = "1:1"
$ws.Range("C29").Value = "اِ۬لْعَٰلَمِينَ"
$ws.Range("D29").Value = "aqaq"
$ws.Range("E29").Value = "Ar"
$ws.Range("F29").Value = "qwqw"
$ws.Range("J29").Value = "24 - 39"
$ws.Range("K29").Value = $false

# Row 30
$ws.Range("B30").Value = "1:1"
$ws.Range("C30").Value = "اِ۬لْعَٰلَمِينَ"
$ws.Range("D30").Value = "aqaq"
$ws.Range("E30").Value = "Ar"
$ws.Range("F30").Value = "qwqw"
$ws.Range("J30").Value = "24 - 39"
$ws.Range("K30").Value = $false

# Row 31
$ws.Range("B31").Value = "1:1"
$ws.Range("C31").Value = "رَبِّ"
$ws.Range("D31").Value = "aqaq"
$ws.Range("E31").Value = "Ar"
$ws.Range("F31").Value = "qwqw"
$ws.Range("J31").Value = "18 - 23"
$ws.Range("K31").Value = $false

# --- 3) Remove the old rows 32-36 (previous rows 30-34) ---------------------
$ws.Range("A32:K36").EntireRow.Delete()

Write-Host "Konduga.xlsx template rows updated"
